# Updates the cryptos list with refreshed prices / 1h volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the column D value
# must be written as text (it looks numeric, and column D cells hold
# plain text such as "234.93" or "1.866.56", not real numbers).
$updates = @(
    @{ Cell = "D2"; Value = "30.196.33"; AsText = $false }
    @{ Cell = "E2"; Value = "  -0.27%  "; AsText = $false }
    @{ Cell = "D3"; Value = "1.864.43"; AsText = $false }
    @{ Cell = "E3"; Value = "  -1.25%  "; AsText = $false }
    @{ Cell = "E4"; Value = "  +0.07%  "; AsText = $false }
    @{ Cell = "D5"; Value = "234.93"; AsText = $true }
    @{ Cell = "E5"; Value = "  -1.48%  "; AsText = $false }
    @{ Cell = "E6"; Value = "  +0.06%  "; AsText = $false }
    @{ Cell = "E7"; Value = "  -0.42%  "; AsText = $false }
    @{ Cell = "D8"; Value = "0.2829"; AsText = $true }
    @{ Cell = "E8"; Value = "  -1.13%  "; AsText = $false }
    @{ Cell = "D9"; Value = "0.06539"; AsText = $true }
    @{ Cell = "E9"; Value = "  -1.05%  "; AsText = $false }
    @{ Cell = "D10"; Value = "21.19"; AsText = $true }
    @{ Cell = "E10"; Value = "  +5.58%  "; AsText = $false }
    @{ Cell = "D11"; Value = "0.07870"; AsText = $true }
    @{ Cell = "E11"; Value = "  +1.00%  "; AsText = $false }
    @{ Cell = "D12"; Value = "97.47"; AsText = $true }
    @{ Cell = "E12"; Value = "  -0.80%  "; AsText = $false }
    @{ Cell = "D13"; Value = "1.867.79"; AsText = $false }
    @{ Cell = "E13"; Value = "  -1.11%  "; AsText = $false }
    @{ Cell = "D14"; Value = "5.099"; AsText = $true }
    @{ Cell = "E14"; Value = "  -0.55%  "; AsText = $false }
    @{ Cell = "D15"; Value = "0.6726"; AsText = $true }
    @{ Cell = "E15"; Value = "  -0.65%  "; AsText = $false }
    @{ Cell = "D16"; Value = "279.68"; AsText = $true }
    @{ Cell = "E16"; Value = "  -2.18%  "; AsText = $false }
    @{ Cell = "D17"; Value = "30.190.00"; AsText = $false }
    @{ Cell = "E17"; Value = "  -0.34%  "; AsText = $false }
    @{ Cell = "D18"; Value = "1.001"; AsText = $true }
    @{ Cell = "E18"; Value = "  +0.06%  "; AsText = $false }
    @{ Cell = "D19"; Value = "5.521"; AsText = $true }
    @{ Cell = "E19"; Value = "  +2.40%  "; AsText = $false }
    @{ Cell = "D20"; Value = "12.66"; AsText = $true }
    @{ Cell = "E20"; Value = "  -0.10%  "; AsText = $false }
    @{ Cell = "D21"; Value = "2.116.53"; AsText = $false }
    @{ Cell = "E21"; Value = "  -1.29%  "; AsText = $false }
    @{ Cell = "D22"; Value = "0.000007277"; AsText = $true }
    @{ Cell = "E22"; Value = "  -0.45%  "; AsText = $false }
    @{ Cell = "D23"; Value = "1.001"; AsText = $true }
    @{ Cell = "E23"; Value = "  +0.05%  "; AsText = $false }
    @{ Cell = "D24"; Value = "6.161"; AsText = $true }
    @{ Cell = "E24"; Value = "  -0.71%  "; AsText = $false }
    @{ Cell = "D25"; Value = "9.194"; AsText = $true }
    @{ Cell = "E25"; Value = "  -2.64%  "; AsText = $false }
    @{ Cell = "D26"; Value = "164.99"; AsText = $true }
    @{ Cell = "E26"; Value = "  -0.57%  "; AsText = $false }
    @{ Cell = "D27"; Value = "19.09"; AsText = $true }
    @{ Cell = "E27"; Value = "  -0.92%  "; AsText = $false }
    @{ Cell = "D28"; Value = "1.920"; AsText = $true }
    @{ Cell = "E28"; Value = "  -3.99%  "; AsText = $false }
    @{ Cell = "D29"; Value = "1.373"; AsText = $true }
    @{ Cell = "E29"; Value = "  -0.35%  "; AsText = $false }
    @{ Cell = "D30"; Value = "0.09700"; AsText = $true }
    @{ Cell = "E30"; Value = "  -0.40%  "; AsText = $false }
    @{ Cell = "D31"; Value = "4.417"; AsText = $true }
    @{ Cell = "E31"; Value = "  -0.68%  "; AsText = $false }
    @{ Cell = "D32"; Value = "1.474"; AsText = $true }
    @{ Cell = "E32"; Value = "  -1.06%  "; AsText = $false }
    @{ Cell = "D33"; Value = "4.096"; AsText = $true }
    @{ Cell = "E33"; Value = "  -1.81%  "; AsText = $false }
    @{ Cell = "D34"; Value = "0.04689"; AsText = $true }
    @{ Cell = "E34"; Value = "  -0.33%  "; AsText = $false }
    @{ Cell = "D35"; Value = "1.117"; AsText = $true }
    @{ Cell = "E35"; Value = "  +1.67%  "; AsText = $false }
    @{ Cell = "D36"; Value = "0.7055"; AsText = $true }
    @{ Cell = "E36"; Value = "  -1.01%  "; AsText = $false }
    @{ Cell = "D37"; Value = "2.726"; AsText = $true }
    @{ Cell = "E37"; Value = "  +0.60%  "; AsText = $false }
    @{ Cell = "D38"; Value = "0.01846"; AsText = $true }
    @{ Cell = "E38"; Value = "  -1.92%  "; AsText = $false }
    @{ Cell = "D39"; Value = "2.537"; AsText = $true }
    @{ Cell = "E39"; Value = "  +0.21%  "; AsText = $false }
    @{ Cell = "D40"; Value = "6.253"; AsText = $true }
    @{ Cell = "E40"; Value = "  -6.95%  "; AsText = $false }
    @{ Cell = "D41"; Value = "73.73"; AsText = $true }
    @{ Cell = "E41"; Value = "  +1.40%  "; AsText = $false }
    @{ Cell = "D42"; Value = "1.941"; AsText = $true }
    @{ Cell = "E42"; Value = "  -1.86%  "; AsText = $false }
    @{ Cell = "D43"; Value = "0.8439"; AsText = $true }
    @{ Cell = "E43"; Value = "  -3.35%  "; AsText = $false }
    @{ Cell = "B44"; Value = "Quant"; AsText = $false }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; AsText = $false }
    @{ Cell = "D44"; Value = "104.00"; AsText = $true }
    @{ Cell = "E44"; Value = "  -0.22%  "; AsText = $false }
    @{ Cell = "B45"; Value = "PaxDollar"; AsText = $false }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; AsText = $false }
    @{ Cell = "D45"; Value = "1.001"; AsText = $true }
    @{ Cell = "E45"; Value = "  +0.05%  "; AsText = $false }
    @{ Cell = "B46"; Value = "TheSandbox"; AsText = $false }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; AsText = $false }
    @{ Cell = "D46"; Value = "0.4164"; AsText = $true }
    @{ Cell = "E46"; Value = "  -1.18%  "; AsText = $false }
    @{ Cell = "D47"; Value = "7.176"; AsText = $true }
    @{ Cell = "E47"; Value = "  -1.37%  "; AsText = $false }
    @{ Cell = "B48"; Value = "EnergySwap"; AsText = $false }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; AsText = $false }
    @{ Cell = "D48"; Value = "9.137"; AsText = $true }
    @{ Cell = "E48"; Value = "  -0.92%  "; AsText = $false }
    @{ Cell = "B49"; Value = "Maker"; AsText = $false }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; AsText = $false }
    @{ Cell = "D49"; Value = "932.93"; AsText = $true }
    @{ Cell = "E49"; Value = "  -5.62%  "; AsText = $false }
    @{ Cell = "D50"; Value = "33.99"; AsText = $true }
    @{ Cell = "E50"; Value = "  -0.51%  "; AsText = $false }
    @{ Cell = "D51"; Value = "0.1124"; AsText = $true }
    @{ Cell = "E51"; Value = "  -3.25%  "; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
